$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "longfish"
$ws.Range("C4").Value = 123

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "longfish"
$ws.Range("C5").Value = 75

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "shortfish"
$ws.Range("C6").Value = 20

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "shortfish"
$ws.Range("C7").Value = 18

$ws.Range("C8").Select()
